$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")
$ws.Range("J2").Value = "Grade 8"
